$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "maxi black is about 35,000x better than uncoated blue PLA"

$ws.Range("C15").Select()
